# Auto-generated Excel COM-interop script to apply scheduled market-data refresh
# to the Siren_Profits leve-profit sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
# Only columns H:N (price/profit columns) change; A:G (leve metadata) are untouched.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 111117860
$ws.Range("I62").Value = 142863420
$ws.Range("J62").Value = 8375
$ws.Range("K62").Value = 142863420
$ws.Range("L62").Value = 8375
$ws.Range("M62").Value = -142862796
$ws.Range("N62").Value = -9623

$ws.Range("H65").Value = 111117860
$ws.Range("I65").Value = 142863420
$ws.Range("J65").Value = 8375
$ws.Range("K65").Value = 714317100
$ws.Range("L65").Value = 41875
$ws.Range("M65").Value = -714313980
$ws.Range("N65").Value = -48115

$ws.Range("H98").Value = 49938.77
$ws.Range("I98").Value = 110880
$ws.Range("J98").Value = 11850.5
$ws.Range("K98").Value = 110880
$ws.Range("L98").Value = 11850.5
$ws.Range("M98").Value = -109382
$ws.Range("N98").Value = -14846.5

$ws.Range("H118").Value = 1878.6666
$ws.Range("I118").Value = 1878.6666
$ws.Range("K118").Value = 5635.9998
$ws.Range("M118").Value = -3978.9998

$ws.Range("H122").Value = 49938.77
$ws.Range("I122").Value = 110880
$ws.Range("J122").Value = 11850.5
$ws.Range("K122").Value = 332640
$ws.Range("L122").Value = 35551.5
$ws.Range("M122").Value = -330190
$ws.Range("N122").Value = -40451.5

$ws.Range("H129").Value = 38463280
$ws.Range("I129").Value = 1863.5
$ws.Range("J129").Value = 50001704
$ws.Range("K129").Value = 5590.5
$ws.Range("L129").Value = 150005112
$ws.Range("M129").Value = -590.5
$ws.Range("N129").Value = -150015112

$ws.Range("H138").Value = 145005
$ws.Range("I138").Value = 824407.8
$ws.Range("J138").Value = 5640.3203
$ws.Range("K138").Value = 2473223.4
$ws.Range("L138").Value = 16920.9609
$ws.Range("M138").Value = -2468083.4
$ws.Range("N138").Value = -27200.9609

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 4875
$ws.Range("I74").Value = 5300
$ws.Range("J74").Value = 4450
$ws.Range("K74").Value = 5300
$ws.Range("L74").Value = 4450
$ws.Range("M74").Value = -4426
$ws.Range("N74").Value = -6198

$ws.Range("H77").Value = 4875
$ws.Range("I77").Value = 5300
$ws.Range("J77").Value = 4450
$ws.Range("K77").Value = 26500
$ws.Range("L77").Value = 22250
$ws.Range("M77").Value = -22132
$ws.Range("N77").Value = -30986

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 4625.5
$ws.Range("I107").Value = 5125.7144
$ws.Range("K107").Value = 5125.7144
$ws.Range("M107").Value = -3205.7144

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2371.3064
$ws.Range("I31").Value = 1554.5555
$ws.Range("J31").Value = 2510
$ws.Range("K31").Value = 1554.5555
$ws.Range("L31").Value = 2510
$ws.Range("M31").Value = -1259.5555
$ws.Range("N31").Value = -3100

$ws.Range("H34").Value = 2371.3064
$ws.Range("I34").Value = 1554.5555
$ws.Range("J34").Value = 2510
$ws.Range("K34").Value = 1554.5555
$ws.Range("L34").Value = 2510
$ws.Range("M34").Value = -1352.5555
$ws.Range("N34").Value = -2914

$ws.Range("H38").Value = 2319
$ws.Range("I38").Value = 2319
$ws.Range("K38").Value = 2319
$ws.Range("M38").Value = -1942

$ws.Range("H39").Value = 2051
$ws.Range("I39").Value = 2051
$ws.Range("K39").Value = 2051
$ws.Range("M39").Value = -1660

$ws.Range("H46").Value = 2319
$ws.Range("I46").Value = 2319
$ws.Range("K46").Value = 2319
$ws.Range("M46").Value = -2108

$ws.Range("H49").Value = 2051
$ws.Range("I49").Value = 2051
$ws.Range("K49").Value = 2051
$ws.Range("M49").Value = -1869

$ws.Range("H57").Value = 52333.332
$ws.Range("J57").Value = 52333.332
$ws.Range("L57").Value = 52333.332
$ws.Range("N57").Value = -53453.332

$ws.Range("H86").Value = 9317.346
$ws.Range("I86").Value = 8050.316
$ws.Range("K86").Value = 8050.316
$ws.Range("M86").Value = -6927.316

$ws.Range("H89").Value = 9317.346
$ws.Range("I89").Value = 8050.316
$ws.Range("K89").Value = 40251.58
$ws.Range("M89").Value = -34635.58

$ws.Range("H130").Value = 79997.5
$ws.Range("J130").Value = 79997.5
$ws.Range("L130").Value = 79997.5
$ws.Range("N130").Value = -90037.5

$ws.Range("H141").Value = 114748.75
$ws.Range("J141").Value = 114748.75
$ws.Range("L141").Value = 114748.75
$ws.Range("N141").Value = -125108.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 358660.2
$ws.Range("J5").Value = 716023
$ws.Range("L5").Value = 2148069
$ws.Range("N5").Value = -2148293

$ws.Range("H38").Value = 1790.7646
$ws.Range("I38").Value = 228.42857
$ws.Range("K38").Value = 685.28571
$ws.Range("M38").Value = -338.28571

$ws.Range("H131").Value = 32261588
$ws.Range("J131").Value = 2092.3635
$ws.Range("L131").Value = 6277.0905
$ws.Range("N131").Value = -16357.0905

$ws.Range("H135").Value = 358660.2
$ws.Range("J135").Value = 716023
$ws.Range("L135").Value = 6444207
$ws.Range("N135").Value = -6449277

$ws.Range("H137").Value = 11177.546
$ws.Range("I137").Value = 6597.6
$ws.Range("K137").Value = 19792.8
$ws.Range("M137").Value = -14692.8

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 23370.572
$ws.Range("I102").Value = 23370.572
$ws.Range("K102").Value = 23370.572
$ws.Range("M102").Value = -21748.572

$ws.Range("H113").Value = 2499.25
$ws.Range("I113").Value = 1498
$ws.Range("K113").Value = 1498
$ws.Range("M113").Value = 672

$ws.Range("H122").Value = 9155.087
$ws.Range("I122").Value = 10917.765
$ws.Range("J122").Value = 4160.8335
$ws.Range("K122").Value = 32753.295
$ws.Range("L122").Value = 12482.5005
$ws.Range("M122").Value = -30303.295
$ws.Range("N122").Value = -17382.5005

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 4913
$ws.Range("I122").Value = 4699.5
$ws.Range("J122").Value = 5340
$ws.Range("K122").Value = 14098.5
$ws.Range("L122").Value = 16020
$ws.Range("M122").Value = -11648.5
$ws.Range("N122").Value = -20920

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 6372.7
$ws.Range("I122").Value = 3498.5
$ws.Range("J122").Value = 9246.9
$ws.Range("K122").Value = 10495.5
$ws.Range("L122").Value = 27740.7
$ws.Range("M122").Value = -8045.5
$ws.Range("N122").Value = -32640.7

$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("N124").ClearContents()

$ws.Range("H126").Value = 20608.584
$ws.Range("I126").Value = 23664.264
$ws.Range("J126").Value = 8997
$ws.Range("K126").Value = 70992.792
$ws.Range("L126").Value = 26991
$ws.Range("M126").Value = -68522.792
$ws.Range("N126").Value = -31931

$ws.Range("H136").Value = 3755.2856
$ws.Range("I136").Value = 3596.7368
$ws.Range("J136").Value = 4090
$ws.Range("K136").Value = 10790.2104
$ws.Range("L136").Value = 12270
$ws.Range("M136").Value = -8240.2104
$ws.Range("N136").Value = -17370
